# daily auto push: 2026-02-05 03:09 UTC
# Insert a new data row for 2026/02/05 10:00 right above the old row 784
# (which starts the 2026/12/29 block), shifting every subsequent row down
# by one. The new row reuses the date/weekday text from row 783 (the
# existing 2026/02/05 05:00 entry) so the date stays a literal text cell
# instead of being re-interpreted as a serial date number, then the
# ranking/time columns are overwritten with the new day's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the existing "2026/02/05" row (A783:D783) so the inserted row
# inherits the same cell types/format (text date + weekday) instead of
# Excel auto-converting the "2026/02/05" string into a date serial.
$ws.Range("A783:D783").Copy()
$ws.Rows("784:784").Insert()

# Overwrite the time-rank (C) and count (D) columns with the new values;
# A784/B784 ("2026/02/05" / "木") came from the copy above and are already
# correct.
$ws.Range("C784").Value = 10
$ws.Range("D784").Value = 201

$excel.CutCopyMode = $false
